$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.129.60'
$ws.Range('E2').Value = '  -2.82%  '
$ws.Range('D3').Value = '3.367.78'
$ws.Range('E3').Value = '  -2.43%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '567.62'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '149.13'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.45%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  +0.26%  '
$ws.Range('E9').Value = '  +1.23%  '
$ws.Range('E10').Value = '  -0.90%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.415'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.92%  '
$ws.Range('D12').Value = '3.943.56'
$ws.Range('E12').Value = '  -2.48%  '
$ws.Range('E13').Value = '  +0.78%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.09'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.90%  '
$ws.Range('D15').Value = '3.372.21'
$ws.Range('E15').Value = '  -2.19%  '
$ws.Range('E16').Value = '  -1.01%  '
$ws.Range('D17').Value = '61.156.77'
$ws.Range('E17').Value = '  -2.88%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.33'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -2.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.48'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -1.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.87'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -2.99%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '375.79'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -3.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '75.34'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.74%  '
$ws.Range('E23').Value = '  -0.33%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').Value = '3.503.74'
$ws.Range('E25').Value = '  -2.53%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000109'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -5.68%  '
$ws.Range('E27').Value = '  -3.79%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.45'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.85%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.19%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.08'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.22%  '
$ws.Range('B31').Value = 'USDe'
$ws.Range('C31').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.07%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.74'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.78%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '22.91'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.69%  '
$ws.Range('E34').Value = '  -3.81%  '
$ws.Range('E35').Value = '  +0.80%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '170.30'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.07%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.55'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -5.23%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.80'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.37%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '29.30'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -7.36%  '
$ws.Range('D40').Value = '3.402.72'
$ws.Range('E40').Value = '  -2.46%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0752'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -3.65%  '
$ws.Range('B42').Value = 'Mantle'
$ws.Range('C42').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.762'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -3.83%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.30'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.18%  '
$ws.Range('B44').Value = 'ONDO'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.15'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.44%  '
$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.61'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -6.20%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '2.490.25'
$ws.Range('E46').Value = '  -2.99%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '22.60'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.25%  '
$ws.Range('E48').Value = '  -3.20%  '
$ws.Range('B49').Value = 'FirstDigitalUSD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.00'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0263'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.67%  '
$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.05'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -9.63%  '
